$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J6 previously held the literal value 900. It now mirrors the shared
# formula in I6 (=I5), extending that shared formula's range to I6:J6.
$ws.Range("J6").Formula = "=I5"

# Update the current selection to match the new active cell/range.
$ws.Range("I6:J6").Select()
